$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles) from column J (rows 2-9) into new column K
$ws.Range("J2:J9").Copy()
$ws.Range("K2:K9").PasteSpecial(-4122)  # xlPasteFormats

# Set the new values for column K
$ws.Range("K3").Value = 2021
$ws.Range("K4").Value = 295
$ws.Range("K5").Value = 163
$ws.Range("K6").Value = 268
$ws.Range("K7").Value = 155
$ws.Range("K8").Value = 27
$ws.Range("K9").Value = 8

# Update the selection to match the target state
$ws.Range("L5").Select()
